# Applies the "Estado de Cuenta" refresh:
#  - Replaces the previous debtor table (3 workers / 4 rows) with the new
#    one (5 workers / 9 periods / 13 rows, rows 16:28).
#  - Moves the signature/footer block from rows 24:25 down to rows 33:34.
#  - Updates the summary fields (Valor Mora total, worker count, period count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Before anything else, capture formatting that needs to move:
#    - row 19's "bottom border" style -> will become row 28's style
#    - rows 24:25's "signature block" style -> will become rows 33:34
# ---------------------------------------------------------------------
$ws.Range("B19:J19").Copy()
$ws.Range("B28:J28").PasteSpecial($xlPasteFormats)

$ws.Range("B24:C24").Copy()
$ws.Range("B33:C33").PasteSpecial($xlPasteFormats)
$ws.Range("H24:J24").Copy()
$ws.Range("H33:J33").PasteSpecial($xlPasteFormats)

$ws.Range("B25:C25").Copy()
$ws.Range("B34:C34").PasteSpecial($xlPasteFormats)
$ws.Range("H25:J25").Copy()
$ws.Range("H34:J34").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Re-create the merges at their new home, drop them from the old rows.
$ws.Range("B33:C33").Merge()
$ws.Range("H33:J33").Merge()
$ws.Range("B34:C34").Merge()
$ws.Range("H34:J34").Merge()

$ws.Range("B24:C24").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("B25:C25").UnMerge()
$ws.Range("H25:J25").UnMerge()

# ---------------------------------------------------------------------
# 2) Rows 19 through 27 become regular (non-last) table rows, so give
#    them the same "middle" formatting already used by rows 16:18.
#    (Row 28 already got its formatting above, from the old row 19.)
# ---------------------------------------------------------------------
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Wipe any leftover values/text on rows 16:28 before refilling - several
# of these rows held old table / footer content that PasteSpecial (format
# only) does not clear.
$ws.Range("B16:J28").ClearContents()

# ---------------------------------------------------------------------
# 3) Fill in the new table: 5 workers, 9 period rows total (16:28).
# ---------------------------------------------------------------------
$table = @(
    @(16, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2206", 25749, 908526),
    @(17, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2205", 36341, 908526),
    @(18, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2204", 36341, 908526),
    @(19, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2203", 36341, 908526),
    @(20, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2202", 36341, 908526),
    @(21, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2201", 36341, 908526),
    @(22, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2112", 36341, 908526),
    @(23, "CC", "73204293",   "ALBERT LADEUTH URUETA",         "2111", 36341, 908526),
    @(24, "CC", "73107141",   "JOHN ALBERT MARTINEZ PALENCIA", "2111", 36341, 908526),
    @(25, "CC", "1047391615", "XAVIER DAVID ORTEGA CHARRY",    "2111", 36341, 908526),
    @(26, "CC", "1047391615", "XAVIER DAVID ORTEGA CHARRY",    "1910", 33125, 908526),
    @(27, "CC", "78751903",   "JUAN CARLOS MUÑOZ MARTINEZ",    "2111", 36341, 908526),
    @(28, "CC", "1047413143", "KEVIN ENRIQUE NAVIA COPETE",    "2111", 36341, 908526)
)

foreach ($row in $table) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: Tipo Doc
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[4]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[5]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[6]   # G: Salario Basico
}

# ---------------------------------------------------------------------
# 4) Re-write the signature / footer block at its new position (33:34).
# ---------------------------------------------------------------------
$ws.Cells.Item(33, 2).Value = "___________________________________"
$ws.Cells.Item(33, 8).Value = "___________________________________"
$ws.Cells.Item(34, 2).Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Cells.Item(34, 8).Value = "FIRMA DEL REPRESENTANTE LEGAL"

# ---------------------------------------------------------------------
# 5) Update the summary header numbers.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 458625   # VALOR MORA total
$ws.Range("C13").Value = 5        # Cant. Trabajadores
$ws.Range("F13").Value = 9        # Cant. Periodos
